# Insert a new data row before the current row 240, shifting all subsequent
# rows (240-332) down by one (to 241-333), and populate the new row with the
# new Alcachofa price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(240).Insert()

$ws.Range("A240").Value2 = 10
$ws.Range("B240").Value2 = 'Vega Modelo de Temuco'
$ws.Range("C240").Value2 = 'La Araucanía'
$ws.Range("D240").Value2 = 45146
$ws.Range("E240").Value2 = 9
$ws.Range("F240").Value2 = 100112013
$ws.Range("G240").Value2 = 'Alcachofa'
$ws.Range("H240").Value2 = 'Madrigal'
$ws.Range("I240").Value2 = 'Primera'
$ws.Range("J240").Value2 = 65
$ws.Range("K240").Value2 = 14000
$ws.Range("L240").Value2 = 14000
$ws.Range("M240").Value2 = 14000
$ws.Range("N240").Value2 = '$/caja 40 unidades'
$ws.Range("O240").Value2 = 'Provincia del Elquí'
$ws.Range("P240").Value2 = 350
$ws.Range("Q240").Value2 = 40
$ws.Range("R240").Value2 = 'Hortaliza'
